$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21
$ws.Range('A21').Value2 = 130872725
$ws.Range('B21').Value2 = 5177
$ws.Range('D21').Value2 = '''LC'
$ws.Range('E21').Value2 = 100526
$ws.Range('F21').Value2 = '''Bronshjon'
$ws.Range('G21').Value2 = '''Callidium coriaceum'
$ws.Range('H21').Value2 = '''Paykull, 1800'
$ws.Range('Q21').Value2 = 570869
$ws.Range('R21').Value2 = 6736590

# Row 22
$ws.Range('A22').Value2 = 130872716
$ws.Range('B22').Value2 = 79243
$ws.Range('D22').Value2 = '''NT'
$ws.Range('E22').Value2 = 6425
$ws.Range('F22').Value2 = '''Garnlav'
$ws.Range('G22').Value2 = '''Alectoria sarmentosa'
$ws.Range('H22').Value2 = '''(Ach.) Ach.'
$ws.Range('Q22').Value2 = 570988
$ws.Range('R22').Value2 = 6736647

# Row 36
$ws.Range('A36').Value2 = 130979080
$ws.Range('B36').Value2 = 57884
$ws.Range('E36').Value2 = 100109
$ws.Range('F36').Value2 = '''Tretåig hackspett'
$ws.Range('G36').Value2 = '''Picoides tridactylus'
$ws.Range('H36').Value2 = '''(Linnaeus, 1758)'
$ws.Range('Q36').Value2 = 571221
$ws.Range('R36').Value2 = 6736517
$ws.Range('S36').Value2 = 1
$ws.Range('AC36').Value2 = '''Äldre ringhack'
$ws.Range('AW36').Value2 = '''Erik Danielsson'
$ws.Range('AX36').Value2 = '''Erik Danielsson'
$ws.Range('Z36').ClearContents()
$ws.Range('AB36').ClearContents()
$ws.Range('AF36').ClearContents()

# Row 37
$ws.Range('A37').Value2 = 130983068
$ws.Range('B37').Value2 = 79243
$ws.Range('E37').Value2 = 6425
$ws.Range('F37').Value2 = '''Garnlav'
$ws.Range('G37').Value2 = '''Alectoria sarmentosa'
$ws.Range('H37').Value2 = '''(Ach.) Ach.'
$ws.Range('Q37').Value2 = 570849
$ws.Range('R37').Value2 = 6736706
$ws.Range('S37').Value2 = 10
$ws.Range('Z37').Value2 = '''09:25'
$ws.Range('AB37').Value2 = '''09:25'
$ws.Range('AW37').Value2 = '''Bo karlstens'
$ws.Range('AX37').Value2 = '''Bo karlstens'
$ws.Range('AC37').ClearContents()
$ws.Range('AF37').NumberFormat = "@"
$ws.Range('AF37').Value2 = ""

# Row 56
$ws.Range('A56').Value2 = 130979100
$ws.Range('B56').Value2 = 79243
$ws.Range('E56').Value2 = 6425
$ws.Range('F56').Value2 = '''Garnlav'
$ws.Range('G56').Value2 = '''Alectoria sarmentosa'
$ws.Range('H56').Value2 = '''(Ach.) Ach.'
$ws.Range('Q56').Value2 = 571473
$ws.Range('R56').Value2 = 6736490

# Row 57
$ws.Range('A57').Value2 = 130979104
$ws.Range('Q57').Value2 = 571129
$ws.Range('R57').Value2 = 6736573
$ws.Range('AF57').NumberFormat = "@"
$ws.Range('AF57').Value2 = ""

# Row 58
$ws.Range('A58').Value2 = 130979094
$ws.Range('Q58').Value2 = 571278
$ws.Range('R58').Value2 = 6736783
$ws.Range('AF58').ClearContents()

# Row 59
$ws.Range('A59').Value2 = 130979086
$ws.Range('B59').Value2 = 91829
$ws.Range('E59').Value2 = 5442
$ws.Range('F59').Value2 = '''Tallticka'
$ws.Range('G59').Value2 = '''Porodaedalea pini'
$ws.Range('H59').Value2 = '''(Brot.) Murrill'
$ws.Range('Q59').Value2 = 571361
$ws.Range('R59').Value2 = 6736509
